$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells AD1:AF1 - copy the formatting of an existing header
# cell (AC1) so the bold/centered/bordered header style is reused, then
# set the header text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Team record (Wins/Losses/Ties) is the same for every player row (2-48).
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 97
    $ws.Cells.Item($r, 31).Value = 65
    $ws.Cells.Item($r, 32).Value = 0
}
